# feat: add 2022-Q1 data
#
# The workbook currently has 4 sheets: 2021-Q2, 2021-Q3, 2021-Q4, 总计 (a
# running summary of "date / holding count / holding value").
#
# This script:
#  1. Turns the old "总计" sheet (index 4) into the new "2022-Q1" sheet,
#     populated with the quarter's per-fund holding detail (same shape as
#     the 2021-Q2 / 2021-Q3 / 2021-Q4 sheets).
#  2. Duplicates that former sheet (pre-edit, so it still carries the
#     original "总计" formatting) into a new sheet placed right after it,
#     renames the duplicate back to "总计", and rewrites its rows to be
#     the summary table with a new leading 2022-Q1 row plus the original
#     2021-Q4 / 2021-Q3 / 2021-Q2 rows shifted down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 0: grab the sheet that is currently "总计" (4th tab) BEFORE any
# renames, and use it as the formatting template for the brand new "总计"
# sheet we will create in step 2.
# ---------------------------------------------------------------------
$sheetQ1 = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------
# Step 1: duplicate it first (while it is still named 总计 + still has the
# old summary content/shape) so the copy inherits all existing formatting,
# then place the duplicate right after the original and rename it.
# ---------------------------------------------------------------------
$sheetQ1.Copy($null, $sheetQ1)
$sheetTotal = $wb.Worksheets.Item(5)

# Rename the original out of the way FIRST (while the duplicate still has
# its auto-generated "总计 (2)"-style name), then rename the duplicate to
# "总计" - doing it in the opposite order would momentarily leave two
# sheets named "总计" and the second rename would be rejected.
$sheetQ1.Name = "2022-Q1"
$sheetTotal.Name = "总计"

# ---------------------------------------------------------------------
# Step 2: rewrite $sheetQ1 ("2022-Q1") with the per-fund holding detail.
# ---------------------------------------------------------------------

# Headers (row 1). B1:D1 already have the "总计-style" header formatting
# (s="2"); extend that same formatting across E1:H1 by copying D1's
# format over, then set all the header text.
$sheetQ1.Range("D1").Copy($sheetQ1.Range("E1:H1")) | Out-Null

$sheetQ1.Cells.Item(1, 2).Value = "基金代码"
$sheetQ1.Cells.Item(1, 3).Value = "基金名称"
$sheetQ1.Cells.Item(1, 4).Value = "基金规模"
$sheetQ1.Cells.Item(1, 5).Value = "股票总仓位"
$sheetQ1.Cells.Item(1, 6).Value = "仓位占比"
$sheetQ1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$sheetQ1.Cells.Item(1, 8).Value = "仓位排名"

function Set-TextCell($ws, $row, $col, $text) {
    # Force the cell to stay text (so numeric-looking strings such as
    # "630010" or "4.93" aren't silently coerced into numbers and don't
    # lose leading zeros / become floats), matching the source data which
    # stores these columns as plain text.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2: 630010 / 华商价值精选混合 / 4.93 / 89.49 / 3.86 / 0.1903 / 5
Set-TextCell $sheetQ1 2 2 "630010"
Set-TextCell $sheetQ1 2 3 "华商价值精选混合"
Set-TextCell $sheetQ1 2 4 "4.93"
Set-TextCell $sheetQ1 2 5 "89.49"
Set-TextCell $sheetQ1 2 6 "3.86"
Set-TextCell $sheetQ1 2 7 "0.1903"
$sheetQ1.Cells.Item(2, 8).Value = 5

# Row 3: 005433 / 申万菱信医药先锋股票 / 2.20 / 90.81 / 4.13 / 0.0909 / 7
Set-TextCell $sheetQ1 3 2 "005433"
Set-TextCell $sheetQ1 3 3 "申万菱信医药先锋股票"
Set-TextCell $sheetQ1 3 4 "2.20"
Set-TextCell $sheetQ1 3 5 "90.81"
Set-TextCell $sheetQ1 3 6 "4.13"
Set-TextCell $sheetQ1 3 7 "0.0909"
$sheetQ1.Cells.Item(3, 8).Value = 7

# Row 4: 630006 / 华商产业升级混合 / 0.98 / 87.95 / 3.79 / 0.0371 / 5
Set-TextCell $sheetQ1 4 2 "630006"
Set-TextCell $sheetQ1 4 3 "华商产业升级混合"
Set-TextCell $sheetQ1 4 4 "0.98"
Set-TextCell $sheetQ1 4 5 "87.95"
Set-TextCell $sheetQ1 4 6 "3.79"
Set-TextCell $sheetQ1 4 7 "0.0371"
$sheetQ1.Cells.Item(4, 8).Value = 5

# ---------------------------------------------------------------------
# Step 3: rewrite $sheetTotal ("总计") with the updated summary table -
# a new 2022-Q1 row on top, the old rows shifted down by one.
# ---------------------------------------------------------------------

# The duplicated sheet only had rows 1-4; row 5 is brand new and needs the
# same row-label formatting (s="2") as column A of the other data rows, so
# copy that format down from A4 before writing the new row 5 value.
$sheetTotal.Range("A4").Copy($sheetTotal.Range("A5")) | Out-Null

$sheetTotal.Cells.Item(2, 1).Value = 0
$sheetTotal.Cells.Item(2, 2).Value = "2022-Q1"
$sheetTotal.Cells.Item(2, 3).Value = 3
$sheetTotal.Cells.Item(2, 4).Value = 0.32

$sheetTotal.Cells.Item(3, 1).Value = 1
$sheetTotal.Cells.Item(3, 2).Value = "2021-Q4"
$sheetTotal.Cells.Item(3, 3).Value = 2
$sheetTotal.Cells.Item(3, 4).Value = 0.29

$sheetTotal.Cells.Item(4, 1).Value = 2
$sheetTotal.Cells.Item(4, 2).Value = "2021-Q3"
$sheetTotal.Cells.Item(4, 3).Value = 11
$sheetTotal.Cells.Item(4, 4).Value = 5.32

$sheetTotal.Cells.Item(5, 1).Value = 3
$sheetTotal.Cells.Item(5, 2).Value = "2021-Q2"
$sheetTotal.Cells.Item(5, 3).Value = 5
$sheetTotal.Cells.Item(5, 4).Value = 1.32
